$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (shifts old N/O/P -> O/P/Q), mirroring the
# "Loan RBI, Variable Instalments" column added to the repayment schedule.
$ws.Columns("N").Insert()

# Match the width Excel assigns the freshly inserted column (closest
# achievable value through the ColumnWidth property's character-width
# rounding).
$ws.Columns("N").ColumnWidth = 9.877604166666666

# Make "Repayment schedule" the active sheet/tab (was "Transactions").
$ws.Activate()

# Update the remembered selection on the now-active sheet.
$ws.Range("M17").Select() | Out-Null
